$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear C2 (cell removed entirely in the fixed output)
$ws.Range("C2").ClearContents()

# Updated values (naive component forecaster bug fix - recalculated forecast vectors)
$ws.Range("E2").Value = 2.829537440099972

$ws.Range("C3").Value = 0.5917823527752386
$ws.Range("E3").Value = -1.985049937500005

$ws.Range("C5").Value = -0.7518797681959066

$ws.Range("C6").Value = -0.5765930039053124

$ws.Range("C7").Value = -0.07642926654479743
$ws.Range("E7").Value = -0.3994003999000184

$ws.Range("C8").Value = 0.9274109147535459
$ws.Range("E8").Value = 2.82953744009995

$ws.Range("E9").Value = -1.194610791899997

$ws.Range("C11").Value = -0.07666472728170559

$ws.Range("E12").Value = -0.3994003999000073

$ws.Range("C13").Value = -0.7283174404323023
$ws.Range("E13").Value = -2.378486270399993

$ws.Range("C14").Value = -0.03096525636256953
$ws.Range("E14").Value = 1.205410808099949

$ws.Range("C15").Value = 1.135350354669407
$ws.Range("E15").Value = 5.718701441600027

$ws.Range("C16").Value = 2.185449115957461

$ws.Range("E19").Value = -3.551690943900021
